$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Simple value updates (rows unaffected by the row-shift) ---
$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# --- From row 9 down, every row shifts because the old "Contact" /
# "No display for ContactDetail" pair (which had been accidentally
# duplicated across rows 10-11) is replaced by a single new
# "Jurisdiction" / "United States of America" row, net removing one row.
# Rewrite rows 9-20 in place with their final values instead of using
# Insert/Delete (which would disturb the existing cell styles), then
# delete the now-unused trailing row 21.

$ws.Range("A9").Value  = "Publisher"
$ws.Range("B9").Value  = "Alvearie Team"

$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "Indicates whether the claim was for a fully insured plan"

$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = $null

$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = $null

$ws.Range("A14").Value = "FHIR Version"
$ws.Range("B14").Value = "4.0.1"

$ws.Range("A15").Value = "Kind"
$ws.Range("B15").Value = "complex-type"

$ws.Range("A16").Value = "Type"
$ws.Range("B16").Value = "Extension"

$ws.Range("A17").Value = "Base Definition"
$ws.Range("B17").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

$ws.Range("A18").Value = "Abstract"
$ws.Range("B18").Value = "false"

$ws.Range("A19").Value = "Derivation"
$ws.Range("B19").Value = "constraint"

$ws.Range("A20").Value = "Context"
$ws.Range("B20").Value = "element:Element"

# Remove the now-superfluous last row (old row 21), shrinking the sheet
# from A1:B21 to A1:B20.
$ws.Rows.Item(21).Delete()
